$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stimulus labels in column B (rows 8-17) to reflect the new
# calibration / VAS scheme agreed after the MRI meeting.
$ws.Range("B8").Value = "pictures\bad2.png"
$ws.Range("B9").Value = "pictures\neutral2.png"
$ws.Range("B10").Value = "VAS2"
$ws.Range("B11").Value = "VAS8"
$ws.Range("B12").Value = "VAS2"
$ws.Range("B13").Value = "VAS8"
$ws.Range("B14").Value = "VAS2"
$ws.Range("B15").Value = "VAS8"
$ws.Range("B16").Value = "VAS2"
$ws.Range("B17").Value = "VAS8"

# Update the selection / scroll position recorded for the sheet view.
$ws.Range("E14").Select() | Out-Null

# Set page setup (paper size + orientation) as agreed for printing.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
